$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.083068580486611
$ws.Range("C2").Value = 0.3889912089759378
$ws.Range("D2").Value = 0.02578659238723624
$ws.Range("E2").Value = 0.07562253673336627
$ws.Range("F2").Value = 1.260833122626849
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 0.7608190033730722
$ws.Range("M2").Value = 0.4319326158130536
$ws.Range("N2").Value = 1.253478116961723

$ws.Range("B3").Value = 0.9674094650519578
$ws.Range("C3").Value = 0.3422684674263223
$ws.Range("D3").Value = 0.02618341187863393
$ws.Range("E3").Value = 0.07133779144054841
$ws.Range("F3").Value = 1.217775299302204
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 0.749140164128498
$ws.Range("M3").Value = 0.3895921539997147
$ws.Range("N3").Value = 1.269049954758621

$ws.Range("B4").Value = 0.8967779099120321
$ws.Range("C4").Value = 0.3136872558378911
$ws.Range("D4").Value = 0.02644378740342823
$ws.Range("E4").Value = 0.06876489150386789
$ws.Range("F4").Value = 1.192292039348587
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 0.7425255544600091
$ws.Range("M4").Value = 0.3638062761704859
$ws.Range("N4").Value = 1.279164933352753

$ws.Range("B5").Value = 0.86808990282276
$ws.Range("C5").Value = 0.3020659771153191
$ws.Range("D5").Value = 0.02655406488560885
$ws.Range("E5").Value = 0.06773075542653828
$ws.Range("F5").Value = 1.182145330265115
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 0.7399689340136817
$ws.Range("M5").Value = 0.3533504694243135
$ws.Range("N5").Value = 1.283425591205251

$ws.Range("B6").Value = 0.8633319701249889
$ws.Range("C6").Value = 0.300137800484066
$ws.Range("D6").Value = 0.02657262739834287
$ws.Range("E6").Value = 0.06755989741226642
$ws.Range("F6").Value = 1.180474786978834
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 0.7395527688929633
$ws.Range("M6").Value = 0.3516174114867638
$ws.Range("N6").Value = 1.284141433036631

$ws.Range("B7").Value = 0.8963906310762582
$ws.Range("C7").Value = 0.3135304241189942
$ws.Range("D7").Value = 0.02644525779485107
$ws.Range("E7").Value = 0.0687508870170106
$ws.Range("F7").Value = 1.192154236666994
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 0.7424905139466631
$ws.Range("M7").Value = 0.3636650556602348
$ws.Range("N7").Value = 1.279221833122982

$ws.Range("B8").Value = 1.043108730473364
$ws.Range("C8").Value = 0.3728584082920747
$ws.Range("D8").Value = 0.02591992372453333
$ws.Range("E8").Value = 0.07413298794923406
$ws.Range("F8").Value = 1.245787428709249
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 0.7566760899402993
$ws.Range("M8").Value = 0.4172891464775361
$ws.Range("N8").Value = 1.258732020899394

$ws.Range("B9").Value = 1.333951536604218
$ws.Range("C9").Value = 0.4901005027218162
$ws.Range("D9").Value = 0.02502398252490856
$ws.Range("E9").Value = 0.08515814239235908
$ws.Range("F9").Value = 1.358635580791073
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 0.7889546083764856
$ws.Range("M9").Value = 0.5241733850616868
$ws.Range("N9").Value = 1.222969436414559

$ws.Range("B10").Value = 1.54968248079058
$ws.Range("C10").Value = 0.5768693767556101
$ws.Range("D10").Value = 0.02444966655351877
$ws.Range("E10").Value = 0.09356162687728187
$ws.Range("F10").Value = 1.446373341710995
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 0.8154579872626613
$ws.Range("M10").Value = 0.6038330271779841
$ws.Range("N10").Value = 1.19941986094809

$ws.Range("B11").Value = 1.64830126175832
$ws.Range("C11").Value = 0.6164983264730495
$ws.Range("D11").Value = 0.02420708247660386
$ws.Range("E11").Value = 0.09745405272480667
$ws.Range("F11").Value = 1.487369310619016
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 0.8281358075262375
$ws.Range("M11").Value = 0.6403357010164967
$ws.Range("N11").Value = 1.189305283742208

$ws.Range("B12").Value = 1.685717092737832
$ws.Range("C12").Value = 0.6315287491458434
$ws.Range("D12").Value = 0.02411794732111616
$ws.Range("E12").Value = 0.09893830691743943
$ws.Range("F12").Value = 1.503051786114696
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 0.8330271019374749
$ws.Range("M12").Value = 0.6541977228741445
$ws.Range("N12").Value = 1.185561801801171

$ws.Range("B13").Value = 1.67765573594545
$ws.Range("C13").Value = 0.6282906013378806
$ws.Range("D13").Value = 0.02413702226814252
$ws.Range("E13").Value = 0.09861818491864227
$ws.Range("F13").Value = 1.499667210743354
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 0.8319696318308871
$ws.Range("M13").Value = 0.6512105245237336
$ws.Range("N13").Value = 1.186364161444374

$ws.Range("B14").Value = 1.651378050666722
$ws.Range("C14").Value = 0.6177344046304256
$ws.Range("D14").Value = 0.02419969440174441
$ws.Range("E14").Value = 0.09757595572350652
$ws.Range("F14").Value = 1.48865633323885
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 0.8285363973644451
$ws.Range("M14").Value = 0.6414753446965165
$ws.Range("N14").Value = 1.188995564026527

$ws.Range("B15").Value = 1.635291514350627
$ws.Range("C15").Value = 0.6112715638606687
$ws.Range("D15").Value = 0.02423843914090895
$ws.Range("E15").Value = 0.09693890648996017
$ws.Range("F15").Value = 1.481932524906512
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 0.826445259144748
$ws.Range("M15").Value = 0.6355174115283404
$ws.Range("N15").Value = 1.190618683370396

$ws.Range("B16").Value = 1.543247354958282
$ws.Range("C16").Value = 0.5742827978181708
$ws.Range("D16").Value = 0.02446589924815434
$ws.Range("E16").Value = 0.09330867140024424
$ws.Range("F16").Value = 1.443716165697822
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 0.8146420608913587
$ws.Range("M16").Value = 0.6014529223013341
$ws.Range("N16").Value = 1.200092966922064

$ws.Range("B17").Value = 1.486905865424319
$ws.Range("C17").Value = 0.5516325179778505
$ws.Range("D17").Value = 0.02461025053717769
$ws.Range("E17").Value = 0.09109966183639528
$ws.Range("F17").Value = 1.420550883784529
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 0.8075611152854378
$ws.Range("M17").Value = 0.580624158679143
$ws.Range("N17").Value = 1.206058825459415

$ws.Range("B18").Value = 1.454544961506883
$ws.Range("C18").Value = 0.5386194157027262
$ws.Range("D18").Value = 0.02469503289430897
$ws.Range("E18").Value = 0.08983563671671391
$ws.Range("F18").Value = 1.407328678409129
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 0.8035467535545138
$ws.Range("M18").Value = 0.5686688887573581
$ws.Range("N18").Value = 1.209546498752104

$ws.Range("B19").Value = 1.443595830568938
$ws.Range("C19").Value = 0.5342158990144981
$ws.Range("D19").Value = 0.02472403899102105
$ws.Range("E19").Value = 0.08940877461518681
$ws.Range("F19").Value = 1.402869288389553
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 0.8021975573637974
$ws.Range("M19").Value = 0.5646252795235824
$ws.Range("N19").Value = 1.210737009267582

$ws.Range("B20").Value = 1.492898819256254
$ws.Range("C20").Value = 0.5540421450093049
$ws.Range("D20").Value = 0.02459470213441861
$ws.Range("E20").Value = 0.0913341360459583
$ws.Range("F20").Value = 1.423006309233614
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 0.8083088405401497
$ws.Range("M20").Value = 0.5828388325993217
$ws.Range("N20").Value = 1.20541792121552

$ws.Range("B21").Value = 1.659094501508889
$ws.Range("C21").Value = 0.6208343601407478
$ws.Range("D21").Value = 0.02418121176135912
$ws.Range("E21").Value = 0.09788180267186419
$ws.Range("F21").Value = 1.49188618255647
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 0.8295423568430209
$ws.Range("M21").Value = 0.6443337286201114
$ws.Range("N21").Value = 1.188220299260223

$ws.Range("B22").Value = 1.768127908675979
$ws.Range("C22").Value = 0.664626213341478
$ws.Range("D22").Value = 0.02392688152776401
$ws.Range("E22").Value = 0.1022211083982043
$ws.Range("F22").Value = 1.537826346321395
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 0.8439476061414837
$ws.Range("M22").Value = 0.6847535241899152
$ws.Range("N22").Value = 1.177486263045218

$ws.Range("B23").Value = 1.709896165592738
$ws.Range("C23").Value = 0.641240540679064
$ws.Range("D23").Value = 0.02406115298251876
$ws.Range("E23").Value = 0.09989955885981772
$ws.Range("F23").Value = 1.513221954022356
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 0.8362105694615423
$ws.Range("M23").Value = 0.6631593567838792
$ws.Range("N23").Value = 1.183168735944157

$ws.Range("B24").Value = 1.490189309746654
$ws.Range("C24").Value = 0.5529527251653121
$ws.Range("D24").Value = 0.02460172598677701
$ws.Range("E24").Value = 0.09122811169943645
$ws.Range("F24").Value = 1.421895913097103
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 0.8079706178409083
$ws.Range("M24").Value = 0.5818375179798778
$ws.Range("N24").Value = 1.205707493972398

$ws.Range("B25").Value = 1.254921411902387
$ws.Range("C25").Value = 0.4582785970090981
$ws.Range("D25").Value = 0.02525177308515758
$ws.Range("E25").Value = 0.08212346921160218
$ws.Range("F25").Value = 1.327270233735646
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 0.7797377847152518
$ws.Range("M25").Value = 0.4950653572952319
$ws.Range("N25").Value = 1.232167932425611
